# "1st changes of mifos to finflux"
# Insert a new (blank) column before column N ("Late") on the
# "Repayment schedule" sheet, pushing the existing Late/Outstanding-
# heading/Outstanding columns one slot to the right, then re-select the
# sheet with the cursor parked at R9 (which becomes the active tab).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column at position 14 (N); existing N..P shift to O..Q.
$ws.Columns.Item(14).Insert()

# The freshly inserted column inherits the width of its left neighbour
# (column M), matching Excel's own insert-column behaviour.
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Make "Repayment schedule" the active sheet/tab with R9 selected.
$ws.Activate()
$ws.Range("R9").Select()
